$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.665.87'
$ws.Range("E2").Value = '  -1.48%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.310.76'
$ws.Range("E3").Value = '  +0.29%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '304.11'
$ws.Range("E5").Value = '  -2.07%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '99.45'
$ws.Range("E6").Value = '  -3.68%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.506'
$ws.Range("E7").Value = '  -4.87%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.503'
$ws.Range("E9").Value = '  -4.68%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.44'
$ws.Range("E10").Value = '  -6.02%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '51.98'
$ws.Range("E11").Value = '  -0.88%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0790'
$ws.Range("E12").Value = '  -2.59%  '
$ws.Range("E13").Value = '  +0.75%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.74'
$ws.Range("E14").Value = '  -3.73%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.673.48'
$ws.Range("E15").Value = '  +0.45%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.64'
$ws.Range("E16").Value = '  +3.71%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.306.22'
$ws.Range("E17").Value = '  -0.01%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.823'
$ws.Range("E18").Value = '  +1.73%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '42.589.24'
$ws.Range("E19").Value = '  -1.42%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0901'
$ws.Range("E20").Value = '  -2.54%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.12'
$ws.Range("E21").Value = '  -0.55%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.52'
$ws.Range("E22").Value = '  -5.22%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '69.16'
$ws.Range("E23").Value = '  +1.61%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '235.03'
$ws.Range("E24").Value = '  -2.98%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.98'
$ws.Range("E25").Value = '  -1.97%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.52'
$ws.Range("E26").Value = '  -3.20%  '
$ws.Range("E27").Value = '  -0.45%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '25.49'
$ws.Range("E28").Value = '  +2.49%  '
$ws.Range("E29").Value = '  +0.01%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.18'
$ws.Range("E30").Value = '  -5.24%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '34.44'
$ws.Range("E31").Value = '  -6.76%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '9.20'
$ws.Range("E32").Value = '  -4.60%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '160.49'
$ws.Range("E33").Value = '  -4.16%  '
$ws.Range("E34").Value = '  -0.03%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.05'
$ws.Range("B36").Value = 'WEMIXToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.45'
$ws.Range("E36").Value = '  -3.12%  '
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.58'
$ws.Range("E37").Value = '  +2.93%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0717'
$ws.Range("E38").Value = '  -3.27%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '16.91'
$ws.Range("E39").Value = '  -7.21%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.88'
$ws.Range("E40").Value = '  -5.38%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.81'
$ws.Range("E41").Value = '  -3.34%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.100'
$ws.Range("E42").Value = '  -5.09%  '
$ws.Range("E43").Value = '  -3.33%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.48'
$ws.Range("E44").Value = '  -5.17%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.994.35'
$ws.Range("E45").Value = '  +0.53%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '18.84'
$ws.Range("E46").Value = '  -0.92%  '
$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0280'
$ws.Range("E47").Value = '  -4.35%  '
$ws.Range("E48").Value = '  +1.88%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.86'
$ws.Range("E49").Value = '  -5.26%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '55.26'
$ws.Range("E50").Value = '  -0.90%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.86'
$ws.Range("E51").Value = '  -2.72%  '
